$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8994777832438388
$ws.Range("D2").Value = 0.03098973583221465
$ws.Range("E2").Value = 0.4989439136740277
$ws.Range("F2").Value = 0.6057878136111157
$ws.Range("G2").Value = 0.4458250048727663
$ws.Range("H2").Value = 0.6048743712228486
$ws.Range("K2").Value = 0.5364747740920848
$ws.Range("L2").Value = 0.09017653330438691
$ws.Range("M2").Value = 0.1716984851608352
$ws.Range("O2").Value = 2.055507917526171

$ws.Range("B3").Value = 0.8793379074288907
$ws.Range("D3").Value = 0.02821835525161021
$ws.Range("E3").Value = 0.503869437435398
$ws.Range("F3").Value = 0.6061101617779272
$ws.Range("G3").Value = 0.4476108319886762
$ws.Range("H3").Value = 0.6091861083516861
$ws.Range("K3").Value = 0.4900690769718494
$ws.Range("L3").Value = 0.08385917601138715
$ws.Range("M3").Value = 0.165806310205344
$ws.Range("O3").Value = 2.068080932139424

$ws.Range("B4").Value = 0.8673459586765375
$ws.Range("D4").Value = 0.02650418037655555
$ws.Range("E4").Value = 0.507072348426215
$ws.Range("F4").Value = 0.6066939203452932
$ws.Range("G4").Value = 0.4490278366037685
$ws.Range("H4").Value = 0.6120999499622783
$ws.Range("K4").Value = 0.461402298598756
$ws.Range("L4").Value = 0.08000511921516562
$ws.Range("M4").Value = 0.1622562861741699
$ws.Range("O4").Value = 2.077029576505197

$ws.Range("B5").Value = 0.8625537122819082
$ws.Range("D5").Value = 0.02580252527819482
$ws.Range("E5").Value = 0.5084225383329728
$ws.Range("F5").Value = 0.6070288643512569
$ws.Range("G5").Value = 0.4496858250490448
$ws.Range("H5").Value = 0.613354412301085
$ws.Range("K5").Value = 0.4496775363109293
$ws.Range("L5").Value = 0.07844090103922241
$ws.Range("M5").Value = 0.1608267925074358
$ws.Range("O5").Value = 2.080985219529765

$ws.Range("B6").Value = 0.8617636905508732
$ws.Range("D6").Value = 0.02568582894961224
$ws.Range("E6").Value = 0.5086494547127538
$ws.Range("F6").Value = 0.6070903446450444
$ws.Range("G6").Value = 0.4497999470849479
$ws.Range("H6").Value = 0.6135667659911732
$ws.Range("K6").Value = 0.4477280862751343
$ws.Range("L6").Value = 0.07818154970352253
$ws.Range("M6").Value = 0.1605904667984355
$ws.Range("O6").Value = 2.081660714186668

$ws.Range("B7").Value = 0.8672809451632304
$ws.Range("D7").Value = 0.02649473017338266
$ws.Range("E7").Value = 0.5070903753632834
$ws.Range("F7").Value = 0.6066980444870111
$ws.Range("G7").Value = 0.4490363843917393
$ws.Range("H7").Value = 0.6121165965310595
$ws.Range("K7").Value = 0.4612443466454579
$ws.Range("L7").Value = 0.07998399781519794
$ws.Range("M7").Value = 0.1622369378396407
$ws.Range("O7").Value = 2.077081672551216

$ws.Range("B8").Value = 0.8924562309124724
$ws.Range("D8").Value = 0.03003678729272252
$ws.Range("E8").Value = 0.5006051996435743
$ws.Range("F8").Value = 0.6058188992556737
$ws.Range("G8").Value = 0.4463742107813999
$ws.Range("H8").Value = 0.6063057963243637
$ws.Range("K8").Value = 0.5205105687489038
$ws.Range("L8").Value = 0.08799321285047768
$ws.Range("M8").Value = 0.1696528634923773
$ws.Range("O8").Value = 2.059588070621132

$ws.Range("B9").Value = 0.9447729837810357
$ws.Range("D9").Value = 0.0368819440857564
$ws.Range("E9").Value = 0.4893021383654741
$ws.Range("F9").Value = 0.6071555778951137
$ws.Range("G9").Value = 0.4436992077807105
$ws.Range("H9").Value = 0.5970224088333396
$ws.Range("K9").Value = 0.6353246722393351
$ws.Range("L9").Value = 0.1038927425874334
$ws.Range("M9").Value = 0.1847288972257708
$ws.Range("O9").Value = 2.035033516118204

$ws.Range("B10").Value = 0.9849859335518261
$ws.Range("D10").Value = 0.04184827217841303
$ws.Range("E10").Value = 0.4818558199128766
$ws.Range("F10").Value = 0.6100037391472739
$ws.Range("G10").Value = 0.4432898313634013
$ws.Range("H10").Value = 0.5914861958614921
$ws.Range("K10").Value = 0.7187882685815907
$ws.Range("L10").Value = 0.1156885951171773
$ws.Range("M10").Value = 0.1961256919213241
$ws.Range("O10").Value = 2.022940252851015

$ws.Range("B11").Value = 1.003661217608425
$ws.Range("D11").Value = 0.04409368302550831
$ws.Range("E11").Value = 0.4786537186033515
$ws.Range("F11").Value = 0.6117047642398461
$ws.Range("G11").Value = 0.4434423750721521
$ws.Range("H11").Value = 0.5892458859586327
$ws.Range("K11").Value = 0.7565582772431014
$ws.Range("L11").Value = 0.1210790254029064
$ws.Range("M11").Value = 0.2013790387553911
$ws.Range("O11").Value = 2.018730782929623

$ws.Range("B12").Value = 1.010787563719617
$ws.Range("D12").Value = 0.0449419450335995
$ws.Range("E12").Value = 0.4774677432314736
$ws.Range("F12").Value = 0.6124071747201612
$ws.Range("G12").Value = 0.4435489145348299
$ws.Range("H12").Value = 0.588437485151303
$ws.Range("K12").Value = 0.7708316142374088
$ws.Range("L12").Value = 0.1231236718714968
$ws.Range("M12").Value = 0.2033781442597373
$ws.Range("O12").Value = 2.017322556552188

$ws.Range("B13").Value = 1.00925036577047
$ws.Range("D13").Value = 0.04475934743064158
$ws.Range("E13").Value = 0.4777219824052921
$ws.Range("F13").Value = 0.6122533069206995
$ws.Range("G13").Value = 0.4435237991630601
$ws.Range("H13").Value = 0.5886098123231136
$ws.Range("K13").Value = 0.767758914747418
$ws.Range("L13").Value = 0.1226831708183909
$ws.Range("M13").Value = 0.2029471683345392
$ws.Range("O13").Value = 2.01761757870014

$ws.Range("B14").Value = 1.004246418799937
$ws.Range("D14").Value = 0.04416351090488035
$ws.Range("E14").Value = 0.4785556150605146
$ws.Range("F14").Value = 0.6117613843159333
$ws.Range("G14").Value = 0.4434501622883147
$ws.Range("H14").Value = 0.5891785777931204
$ws.Range("K14").Value = 0.7577331463850498
$ws.Range("L14").Value = 0.1212471720174193
$ws.Range("M14").Value = 0.2015433111593197
$ws.Range("O14").Value = 2.018611203160503

$ws.Range("B15").Value = 1.001188428034709
$ws.Range("D15").Value = 0.04379827867965957
$ws.Range("E15").Value = 0.4790697005700721
$ws.Range("F15").Value = 0.6114676551387817
$ws.Range("G15").Value = 0.4434114111978289
$ws.Range("H15").Value = 0.589532165307844
$ws.Range("K15").Value = 0.7515882262649711
$ws.Range("L15").Value = 0.1203680218832091
$ws.Range("M15").Value = 0.2006846775430873
$ws.Range("O15").Value = 2.019244025928771

$ws.Range("B16").Value = 0.9837731047018963
$ws.Range("D16").Value = 0.04170124805074238
$ws.Range("E16").Value = 0.482068808498501
$ws.Range("F16").Value = 0.609900728444245
$ws.Range("G16").Value = 0.4432866843397534
$ws.Range("H16").Value = 0.5916381988387371
$ws.Range("K16").Value = 0.7163158506645573
$ws.Range("L16").Value = 0.1153368004557649
$ws.Range("M16").Value = 0.1957837485274752
$ws.Range("O16").Value = 2.023241350334729

$ws.Range("B17").Value = 0.9731868769395646
$ws.Range("D17").Value = 0.04041122386211526
$ws.Range("E17").Value = 0.4839560793973963
$ws.Range("F17").Value = 0.6090432810040483
$ws.Range("G17").Value = 0.443296978691734
$ws.Range("H17").Value = 0.5930013920643376
$ws.Range("K17").Value = 0.6946260553768866
$ws.Range("L17").Value = 0.1122564922895606
$ws.Range("M17").Value = 0.1927947393755645
$ws.Range("O17").Value = 2.026024478842828

$ws.Range("B18").Value = 0.9671339665459868
$ws.Range("D18").Value = 0.03966793999403251
$ws.Range("E18").Value = 0.4850590288196361
$ws.Range("F18").Value = 0.6085882579604487
$ws.Range("G18").Value = 0.4433347829911156
$ws.Range("H18").Value = 0.5938116471322843
$ws.Range("K18").Value = 0.6821320740479848
$ws.Range("L18").Value = 0.1104870870742616
$ws.Range("M18").Value = 0.1910820321099322
$ws.Range("O18").Value = 2.027746849984595

$ws.Range("B19").Value = 0.9650907559351936
$ws.Range("D19").Value = 0.03941605553651328
$ws.Range("E19").Value = 0.4854354656063267
$ws.Range("F19").Value = 0.6084407493723987
$ws.Range("G19").Value = 0.4433530570208575
$ws.Range("H19").Value = 0.5940904836569416
$ws.Range("K19").Value = 0.6778986613803681
$ws.Range("L19").Value = 0.1098883962765456
$ws.Range("M19").Value = 0.190503257750521
$ws.Range("O19").Value = 2.028350896851165

$ws.Range("B20").Value = 0.9743100761981509
$ws.Range("D20").Value = 0.04054868358478103
$ws.Range("E20").Value = 0.4837533714629449
$ws.Range("F20").Value = 0.6091306087043478
$ws.Range("G20").Value = 0.4432925826064107
$ws.Range("H20").Value = 0.5928535683330693
$ws.Range("K20").Value = 0.696936900904717
$ws.Range("L20").Value = 0.1125841584041183
$ws.Range("M20").Value = 0.1931122535657437
$ws.Range("O20").Value = 2.025715625681954

$ws.Range("B21").Value = 1.005714726645778
$ws.Range("D21").Value = 0.04433857786825968
$ws.Range("E21").Value = 0.4783100356707193
$ws.Range("F21").Value = 0.6119042927013751
$ws.Range("G21").Value = 0.4434704670299396
$ws.Range("H21").Value = 0.5890104335173163
$ws.Range("K21").Value = 0.7606787609516346
$ws.Range("L21").Value = 0.1216688680450915
$ws.Range("M21").Value = 0.2019553937217111
$ws.Range("O21").Value = 2.01831430842185

$ws.Range("B22").Value = 1.026556440024109
$ws.Range("D22").Value = 0.04680366205708708
$ws.Range("E22").Value = 0.4749074701627505
$ws.Range("F22").Value = 0.6140566806289414
$ws.Range("G22").Value = 0.4438710464242064
$ws.Range("H22").Value = 0.586731597709317
$ws.Range("K22").Value = 0.8021663004777224
$ws.Range("L22").Value = 0.127626063100152
$ws.Range("M22").Value = 0.2077918255138442
$ws.Range("O22").Value = 2.014560173659561

$ws.Range("B23").Value = 1.015404001967283
$ws.Range("D23").Value = 0.04548909677380664
$ws.Range("E23").Value = 0.4767093202040842
$ws.Range("F23").Value = 0.612876843388598
$ws.Range("G23").Value = 0.4436312148721555
$ws.Range("H23").Value = 0.5879265605623942
$ws.Range("K23").Value = 0.7800395951011012
$ws.Range("L23").Value = 0.1244448183798283
$ws.Range("M23").Value = 0.2046716456306399
$ws.Range("O23").Value = 2.016464711802428

$ws.Range("B24").Value = 0.973802174193537
$ws.Range("D24").Value = 0.04048654312376954
$ws.Range("E24").Value = 0.4838449598594599
$ws.Range("F24").Value = 0.609091009685514
$ws.Range("G24").Value = 0.4432944707588078
$ws.Range("H24").Value = 0.5929203167772386
$ws.Range("K24").Value = 0.6958922430469556
$ws.Range("L24").Value = 0.112436015858421
$ws.Range("M24").Value = 0.1929686875927885
$ws.Range("O24").Value = 2.025854877202391

$ws.Range("B25").Value = 0.9303065176754899
$ws.Range("D25").Value = 0.03504107969841641
$ws.Range("E25").Value = 0.4922089341791551
$ws.Range("F25").Value = 0.6064662750720444
$ws.Range("G25").Value = 0.4441499204355353
$ws.Range("H25").Value = 0.5993080607749874
$ws.Range("K25").Value = 0.6044184191307522
$ws.Range("L25").Value = 0.09957114552686619
$ws.Range("M25").Value = 0.1805938161819078
$ws.Range("O25").Value = 2.019244025928771
